$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the climate3 label: remove embedded line break (A10 holds this string)
$ws.Range("A10").Value = "climate3: Ban the sale of new combustion-engine cars by 2030"

# Overwrite column B values with corrected figures (right-size past figures)
$ws.Range("B2").Value = 13.4380211658986
$ws.Range("B3").Value = 23.0388403171059
$ws.Range("B4").Value = 21.2287024992306
$ws.Range("B5").Value = 27.7651677414388
$ws.Range("B6").Value = 10.1906132579199
$ws.Range("B7").Value = 13.3797440916754
$ws.Range("B8").Value = 13.7042389596488
$ws.Range("B9").Value = 19.7507435781098
$ws.Range("B10").Value = 10.5530473217693
$ws.Range("B11").Value = 13.5810946595551
$ws.Range("B12").Value = 19.0358810776434
$ws.Range("B13").Value = 15.4001299565062
$ws.Range("B14").Value = 20.630826845836
$ws.Range("B15").Value = 14.646697373968
$ws.Range("B16").Value = 8.82406787750107
